$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregado nuevo caso de prueba: Validar cupón inválido
$ws.Range("A7").Value = "TC06"
$ws.Range("B7").Value = "Validar Cupón inválido"
$ws.Range("E7").Value = "Mostrar mensaje de error claro "

$ws.Range("E7").Select()
